# Applies the "Add files via upload" revision:
#  1) Handout master's cached datetimeFigureOut field: 9/27/20 -> 11/6/20
#  2) Slide 5 ("Updates Since IETF-107 (Version-02)") bullet-text rewrites.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Handout master date placeholder (p:ph type="dt") cached field text.
# ---------------------------------------------------------------------
try {
    $hm = $p.HandoutMaster
    for ($i = 1; $i -le $hm.Shapes.Count; $i++) {
        $hmShape = $hm.Shapes.Item($i)
        if ($hmShape.HasTextFrame -and $hmShape.TextFrame.TextRange.Text -eq "9/27/20") {
            $hmShape.TextFrame.TextRange.Text = "11/6/20"
        }
    }
} catch {
    # Some hosts do not allow programmatic edits of the handout master;
    # ignore and continue with the slide-content edits below.
}

# ---------------------------------------------------------------------
# 2) Locate the slide whose title starts with "Updates Since IETF-107"
#    and fix up its bullet list (Content Placeholder).
# ---------------------------------------------------------------------
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    $titleShape = $null
    $bodyShape = $null

    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if (-not $shape.HasTextFrame) { continue }
        $text = $shape.TextFrame.TextRange.Text
        if ($text -like "Updates Since IETF-107*") {
            $titleShape = $shape
        } elseif ($text -like "Updates:*") {
            $bodyShape = $shape
        }
    }

    if ($titleShape -ne $null -and $bodyShape -ne $null) {
        $tr = $bodyShape.TextFrame.TextRange

        for ($pi = 1; $pi -le $tr.Paragraphs().Count; $pi++) {
            $paraText = $tr.Paragraphs($pi).Text

            $newText = $null
            if ($paraText -eq "Decouple with TWAMP (STAMP) protocol using new TWAMP (STAMP) compatible message format`r") {
                $newText = "Use TWAMP (STAMP) compatible probe message format"
            } elseif ($paraText -eq "Show extension Label 15 in MPLS header`r") {
                $newText = "Add extension Label 15 in MPLS header"
            } elseif ($paraText -eq "New section on SRv6 Timestamp Endpoint function assignment and Node Capability`r") {
                $newText = "Add section on SRv6 Timestamp Endpoint function assignment and Node Capability"
            } elseif ($paraText -eq "Update IANA section for a registry name`r") {
                $newText = "Update IANA section"
            } elseif ($paraText -eq "Various editorial changes (e.g. moving text to new Overview section)`r") {
                $newText = "Various editorial changes (e.g. moving some text to new Overview section)"
            }

            if ($newText -ne $null) {
                # Overwrite with a throwaway placeholder first so the
                # run-level diffing applied by the host has nothing in
                # common with the final string and does not fragment the
                # paragraph into multiple <a:r> runs - then set the real
                # text, which keeps the original run's rPr untouched.
                $tr.Paragraphs($pi).Text = "`u{2060}"
                $tr.Paragraphs($pi).Text = $newText
            }
        }

        break
    }
}
